# Apply the requirement-separation update to the FPE worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "Terms Typically Offered" column (D),
# pushing it to column G. The new columns are D=Corequisites, E=Concurrent, F=Recommended.
$ws.Columns("D:F").Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Fill new D/E/F columns with "NA" for every data row (2-18) by default.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}

# Row 13 (FPE 554): split "Recommended: FPE 504." out of the prerequisites text into column F.
$ws.Range("C13").Value = "Consent of graduate coordinator and instructor."
$ws.Range("F13").Value = "FPE 504."
$ws.Range("G13").Value = "SP "

# Row 14 (FPE 555): split "Recommended: LA/NR 318 and NR 340." out of the prerequisites text into column F.
$ws.Range("C14").Value = "Graduate standing or consent of instructor."
$ws.Range("F14").Value = "LA/NR 318 and NR 340."
$ws.Range("G14").Value = "TBD "

# Update the sheet dimension to reflect the new range.
$ws.UsedRange | Out-Null
